# Auto-generated Excel COM-interop script to apply translation updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: Delete obsolete rows (descending order) ---
$ws.Rows(124).Delete() | Out-Null
$ws.Rows(12).Delete() | Out-Null
$ws.Rows(11).Delete() | Out-Null

# --- Step 2: Insert new rows at their final target positions (ascending order) ---
$ws.Rows(38).Insert() | Out-Null
$ws.Cells.Item(38, 1).Value = 'Contains names of organisms before and after mapping.'
$ws.Cells.Item(38, 2).Value = 'Contient les noms des organismes avant et après cartorgraphie.'

$ws.Rows(50).Insert() | Out-Null
$ws.Cells.Item(50, 1).Value = 'Download Lab Log (.xlsx)'
$ws.Cells.Item(50, 2).Value = 'Télécharger un log du traitement des donnés de lab (.xlsx)'

$ws.Rows(110).Insert() | Out-Null
$ws.Cells.Item(110, 1).Value = 'Remove ''Not Cultured'' specimens'
$ws.Cells.Item(110, 2).Value = 'Supprimer les spécimens ''Not Cultured'''

# --- Step 3: Update cell values (fill in translations / text edits) ---
$ws.Cells.Item(10, 2).Value = '(Pour vous déconnecter, fermer l''app.)'
$ws.Cells.Item(20, 2).Value = 'Tous les enregistrements valides ont une ID ACORN.'
$ws.Cells.Item(21, 2).Value = 'AMR'
$ws.Cells.Item(35, 2).Value = 'Co-résistances'
$ws.Cells.Item(47, 2).Value = 'Annuler'
$ws.Cells.Item(65, 2).Value = 'Générer un fichier .acorn depuis les données cliniques et de lab'
$ws.Cells.Item(68, 2).Value = 'Obtenir la dernière version stable'
$ws.Cells.Item(70, 1).Value = 'HAI point prevalence by '
$ws.Cells.Item(70, 2).Value = 'TBT'
$ws.Cells.Item(72, 2).Value = 'Les barres horizontales indiquent la taille d''un ensemble de résultats SR tandis que les barres verticales indiquent le nombre d''isolats résistants pour l''antibiotique correspondant.'
$ws.Cells.Item(73, 2).Value = 'Information sur le fichier .acorn chargé.'
$ws.Cells.Item(76, 2).Value = 'Problème détecté avec les données REDCap. Merci de contacter l''équipe ACORN. Jusqu''à résolution, seuls les fichiers .acorn existants peuvent être utilisés.'
$ws.Cells.Item(84, 2).Value = 'Charger le fichier .acorn depuis le nuage'
$ws.Cells.Item(85, 2).Value = 'Charger le fichier .acorn localement'
$ws.Cells.Item(100, 2).Value = 'Seuls les isolats qui ont été testés contre tous les médicaments sont inclus dans le graphique.'
$ws.Cells.Item(125, 2).Value = 'Montrer les combinaisons d''antibiotiques.'
$ws.Cells.Item(127, 2).Value = 'Evaluation SIR'
$ws.Cells.Item(133, 2).Value = 'Certains enregistrements ont un identifiant ACORN manquant. Ces enregistrements ont été supprimés.'
$ws.Cells.Item(135, 2).Value = 'Spécimens'
$ws.Cells.Item(144, 2).Value = 'Susceptible & Intermédiaire sont toujours combinés dans cette visualisation des co-résistances.'
$ws.Cells.Item(146, 2).Value = 'Les « identifiants de patient » suivants sont des cas atypiques (un HCAI/CAI avec HAI précoce mais sans chevauchement) :'
$ws.Cells.Item(148, 2).Value = 'Le jeu de données REDCap est vide/au mauvais format. Veuillez contacter l''assistance ACORN.'
$ws.Cells.Item(149, 2).Value = 'Le jeu de données REDCap est au bon format.'
$ws.Cells.Item(150, 2).Value = 'Il y a des suivis à J28 effectués avant la date prévue à J28.'
$ws.Cells.Item(151, 2).Value = 'Il existe plusieurs F02 avec un ID ACORN, une date d''admission et une date d''enrôlement identiques.'
$ws.Cells.Item(152, 2).Value = 'Il n''y a pas de cas atypique (un HCAI/CAI avec HAI précoce mais pas de chevauchement).'
$ws.Cells.Item(153, 2).Value = 'Il n''y a pas de suivi à J28 effectué avant la date prévue à J28.'
$ws.Cells.Item(155, 2).Value = 'Il n''existe pas de F02 avec un ID ACORN, une date d''admission et une date d''enrôlement identiques.'
$ws.Cells.Item(167, 2).Value = 'Updated Charlson Comorbidity Index (uCCI)'
$ws.Cells.Item(172, 2).Value = 'Le dictionnaire des données de lab ne peut pas être téléchargé. Merci de contacter l''équipe ACORN.'
$ws.Cells.Item(177, 2).Value = 'Vous utilisez le tableau de bord ACORN'
$ws.Cells.Item(178, 2).Value = 'Vous pouvez vérifier ici s''il s''agit de la dernière version de production.'
$ws.Cells.Item(179, 2).Value = 'Votre tableau de bord ACORN est à jour'
